$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the curated dimension/measure metadata values (row 2)
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("H2").Value = "iaest-measure:tipo-estudios"

# Update the dim/medida classification (row 3)
$ws.Range("E3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Update the value type metadata (row 4)
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"

# Remove the now-obsolete mapping file row entirely
$ws.Rows("5").Delete()
